$wb = $excel.ActiveWorkbook

# @@ -6298,25 +6298,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 14746420
$ws.Range("I113").Value = 8335957.5
$ws.Range("J113").Value = 25003160
$ws.Range("K113").Value = 8335957.5
$ws.Range("L113").Value = 25003160
$ws.Range("M113").Value = -8332703.5
$ws.Range("N113").Value = -25009668

# @@ -6396,22 +6396,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 899.26666
$ws.Range("I115").Value = 499
$ws.Range("K115").Value = 1497
$ws.Range("M115").Value = 70

# @@ -6546,25 +6546,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 853.9524
$ws.Range("I118").Value = 475.53333
$ws.Range("J118").Value = 1800
$ws.Range("K118").Value = 1426.59999
$ws.Range("L118").Value = 5400
$ws.Range("M118").Value = 230.4000100000001
$ws.Range("N118").Value = -8714

# @@ -7100,25 +7100,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1001.13336
$ws.Range("I129").Value = 592.4286
$ws.Range("J129").Value = 1125.5217
$ws.Range("K129").Value = 1777.2858
$ws.Range("L129").Value = 3376.5651
$ws.Range("M129").Value = 3222.7142
$ws.Range("N129").Value = -13376.5651

# @@ -9310,22 +9310,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6664.459
$ws.Range("I32").Value = 4787.463
$ws.Range("K32").Value = 4787.463
$ws.Range("M32").Value = -4500.463

# @@ -9953,25 +9953,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 13205.75
$ws.Range("I45").Value = 13205.75
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 13205.75
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -12828.75
$ws.Range("N45").ClearContents()

# @@ -13684,25 +13681,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 733995
$ws.Range("I122").Value = 885556.4399999999
$ws.Range("J122").Value = 1448.3334
$ws.Range("K122").Value = 2656669.32
$ws.Range("L122").Value = 4345.0002
$ws.Range("M122").Value = -2654219.32
$ws.Range("N122").Value = -9245.0002

# @@ -14168,22 +14165,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3693.24
$ws.Range("I132").Value = 1965.4736
$ws.Range("K132").Value = 5896.4208
$ws.Range("M132").Value = -3366.4208

# @@ -19787,25 +19784,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10588.417
$ws.Range("I105").Value = 15559.214
$ws.Range("J105").Value = 3629.3
$ws.Range("K105").Value = 15559.214
$ws.Range("L105").Value = 3629.3
$ws.Range("M105").Value = -13812.214
$ws.Range("N105").Value = -7123.3

# @@ -21205,25 +21202,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2116.9412
$ws.Range("I134").Value = 2116.9412
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6350.823600000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3815.823600000001
$ws.Range("N134").ClearContents()

# @@ -22383,22 +22377,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1989.2941
$ws.Range("I16").Value = 1678.3334
$ws.Range("K16").Value = 1678.3334
$ws.Range("M16").Value = -1391.3334

# @@ -26765,22 +26759,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2370.818
$ws.Range("I105").Value = 2309.875
$ws.Range("K105").Value = 2309.875
$ws.Range("M105").Value = -562.875

# @@ -27160,22 +27154,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1989.2941
$ws.Range("I113").Value = 1678.3334
$ws.Range("K113").Value = 1678.3334
$ws.Range("M113").Value = 491.6666

# @@ -27595,22 +27589,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2780318.8
$ws.Range("I122").Value = 4630029
$ws.Range("K122").Value = 13890087
$ws.Range("M122").Value = -13887637

# @@ -28073,22 +28067,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2912
$ws.Range("I132").Value = 2328.6667
$ws.Range("K132").Value = 6986.000100000001
$ws.Range("M132").Value = -4456.000100000001

# @@ -28174,22 +28168,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2510.4878
$ws.Range("I134").Value = 2405.6843
$ws.Range("K134").Value = 7217.0529
$ws.Range("M134").Value = -4682.0529

# @@ -28718,25 +28712,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7594.346
$ws.Range("I3").Value = 10651.667
$ws.Range("J3").Value = 4973.7856
$ws.Range("K3").Value = 31955.001
$ws.Range("L3").Value = 14921.3568
$ws.Range("M3").Value = -31843.001
$ws.Range("N3").Value = -15145.3568

# @@ -28822,25 +28816,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 223168.38
$ws.Range("I5").Value = 499
$ws.Range("J5").Value = 346020.44
$ws.Range("K5").Value = 1497
$ws.Range("L5").Value = 1038061.32
$ws.Range("M5").Value = -1385
$ws.Range("N5").Value = -1038285.32

# @@ -34237,25 +34231,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1000553.8
$ws.Range("J113").Value = 1250521.2
$ws.Range("L113").Value = 3751563.6
$ws.Range("N113").Value = -3755903.6

# @@ -34702,25 +34696,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1973.5146
$ws.Range("J122").Value = 2231.1553
$ws.Range("L122").Value = 20080.3977
$ws.Range("N122").Value = -24980.3977

# @@ -35372,25 +35366,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 223168.38
$ws.Range("I135").Value = 499
$ws.Range("J135").Value = 346020.44
$ws.Range("K135").Value = 4491
$ws.Range("L135").Value = 3114183.96
$ws.Range("M135").Value = -1956
$ws.Range("N135").Value = -3119253.96

# @@ -35476,25 +35470,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6206.1763
$ws.Range("I137").Value = 5011.1113
$ws.Range("J137").Value = 7550.625
$ws.Range("K137").Value = 15033.3339
$ws.Range("L137").Value = 22651.875
$ws.Range("M137").Value = -9933.333899999998
$ws.Range("N137").Value = -32851.875

# @@ -41891,25 +41885,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4725.3486
$ws.Range("I126").Value = 10081.167
$ws.Range("J126").Value = 2652.1292
$ws.Range("K126").Value = 30243.501
$ws.Range("L126").Value = 7956.3876
$ws.Range("M126").Value = -27773.501
$ws.Range("N126").Value = -12896.3876

# @@ -42179,25 +42173,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2226.2
$ws.Range("I132").Value = 1703.3334
$ws.Range("J132").Value = 3010.5
$ws.Range("K132").Value = 5110.0002
$ws.Range("L132").Value = 9031.5
$ws.Range("M132").Value = -2580.0002
$ws.Range("N132").Value = -14091.5

# @@ -43008,22 +43002,22 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49229.047
$ws.Range("I7").Value = 54095.26
$ws.Range("K7").Value = 54095.26
$ws.Range("M7").Value = -53983.26

# @@ -44634,25 +44628,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 41669010
$ws.Range("I40").Value = 58824964
$ws.Range("J40").Value = 4549.2856
$ws.Range("K40").Value = 58824964
$ws.Range("L40").Value = 4549.2856
$ws.Range("M40").Value = -58824828
$ws.Range("N40").Value = -4821.2856

# @@ -45651,25 +45645,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1798.1111
$ws.Range("I61").Value = 1769.7142
$ws.Range("J61").Value = 1897.5
$ws.Range("K61").Value = 1769.7142
$ws.Range("L61").Value = 1897.5
$ws.Range("M61").Value = -1567.7142
$ws.Range("N61").Value = -2301.5

# @@ -48154,25 +48148,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1798.1111
$ws.Range("I113").Value = 1769.7142
$ws.Range("J113").Value = 1897.5
$ws.Range("K113").Value = 1769.7142
$ws.Range("L113").Value = 1897.5
$ws.Range("M113").Value = 400.2858000000001
$ws.Range("N113").Value = -6237.5

# @@ -48580,25 +48574,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4288752
$ws.Range("I122").Value = 7938998
$ws.Range("J122").Value = 1003530.8
$ws.Range("K122").Value = 23816994
$ws.Range("L122").Value = 3010592.4
$ws.Range("M122").Value = -23814544
$ws.Range("N122").Value = -3015492.4

# @@ -48776,22 +48770,22 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 49229.047
$ws.Range("I126").Value = 54095.26
$ws.Range("K126").Value = 162285.78
$ws.Range("M126").Value = -159815.78

# @@ -49064,25 +49058,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 37040628
$ws.Range("I132").Value = 111112950
$ws.Range("J132").Value = 4467.1665
$ws.Range("K132").Value = 333338850
$ws.Range("L132").Value = 13401.4995
$ws.Range("M132").Value = -333336320
$ws.Range("N132").Value = -18461.4995

# @@ -54784,25 +54778,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 250000820
$ws.Range("I107").Value = 250000820
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 750002460
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -750000540
$ws.Range("N107").ClearContents()

# @@ -55510,22 +55501,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1009.8
$ws.Range("I122").Value = 1009.8
$ws.Range("K122").Value = 3029.4
$ws.Range("M122").Value = -579.3999999999996

# @@ -55703,25 +55694,25 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 711.3103599999999
$ws.Range("I126").Value = 434.65
$ws.Range("J126").Value = 1326.1111
$ws.Range("K126").Value = 1303.95
$ws.Range("L126").Value = 3978.3333
$ws.Range("M126").Value = 1166.05
